# chore: update Sheets via scheduled runner
# Applies the scheduled price-refresh update to the Leve profit tables
# across the ALC / ARM / BSM / CRP / CUL / LTW / WVR sheets.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $val = $Values[$col]
        $cell = $ws.Range("$col$Row")
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# ---- ALC ----
Set-Row "ALC" 137 @{ H = 1458.6198; I = 1074.8628; K = 3224.588400000001; M = -674.5884000000005 }

# ---- ARM ----
Set-Row "ARM" 32  @{ H = 18898.586; I = 19854.8; K = 19854.8; M = -19567.8 }
Set-Row "ARM" 45  @{ H = 1332.017; I = 1254.7959; K = 1254.7959; M = -877.7959000000001 }
Set-Row "ARM" 61  @{ H = 7211.0576; I = 3892.0952; J = 21150.7; K = 3892.0952; L = 21150.7; M = -3680.0952; N = -21574.7 }
Set-Row "ARM" 74  @{ H = 4220.1313; I = 1798.6207; J = 12022.777; K = 1798.6207; L = 12022.777; M = -924.6206999999999; N = -13770.777 }
Set-Row "ARM" 77  @{ H = 4220.1313; I = 1798.6207; J = 12022.777; K = 8993.103499999999; L = 60113.885; M = -4625.103499999999; N = -68849.88500000001 }
Set-Row "ARM" 102 @{ H = 2649101.8; I = 4632753; J = 4233.1665; K = 4632753; L = 4233.1665; M = -4631131; N = -7477.1665 }
Set-Row "ARM" 136 @{ H = 7211.0576; I = 3892.0952; J = 21150.7; K = 11676.2856; L = 63452.10000000001; M = -9126.285600000001; N = -68552.10000000001 }

# ---- BSM ----
Set-Row "BSM" 22  @{ H = 334.58334; I = 319.54544; K = 319.54544; M = -146.54544 }
Set-Row "BSM" 86  @{ H = 1630; I = 1629.6875; J = 1633.3334; K = 1629.6875; L = 1633.3334; M = -506.6875; N = -3879.3334 }
Set-Row "BSM" 89  @{ H = 1630; I = 1629.6875; J = 1633.3334; K = 8148.4375; L = 8166.666999999999; M = -2532.4375; N = -19398.667 }
Set-Row "BSM" 105 @{ H = 923740.7; I = 1306486.8; J = 5150; K = 1306486.8; L = 5150; M = -1304739.8; N = -8644 }

# ---- CRP ----
Set-Row "CRP" 31  @{ H = 1790.2; I = 1228.2037; J = 4549.091; K = 1228.2037; L = 4549.091; M = -933.2037; N = -5139.091 }
Set-Row "CRP" 34  @{ H = 1790.2; I = 1228.2037; J = 4549.091; K = 1228.2037; L = 4549.091; M = -1026.2037; N = -4953.091 }
Set-Row "CRP" 62  @{ H = 2687.5; I = 2894.111; J = 2067.6667; K = 2894.111; L = 2067.6667; M = -2270.111; N = -3315.6667 }
Set-Row "CRP" 65  @{ H = 2687.5; I = 2894.111; J = 2067.6667; K = 14470.555; L = 10338.3335; M = -11350.555; N = -16578.3335 }
Set-Row "CRP" 68  @{ H = 38967.145; J = 38967.145; L = 38967.145; N = -40465.145 }
Set-Row "CRP" 71  @{ H = 38967.145; J = 38967.145; L = 116901.435; N = -124389.435 }
Set-Row "CRP" 88  @{ H = 40343; I = 0; J = 40343; K = 0; L = 40343; M = $null; N = -41155 }
Set-Row "CRP" 91  @{ H = 40343; I = 0; J = 40343; K = 0; L = 40343; M = $null; N = -43151 }
Set-Row "CRP" 132 @{ H = 2987.3108; I = 3310.0193; K = 9930.0579; M = -7400.0579 }

# ---- CUL ----
Set-Row "CUL" 44  @{ H = 823.75; I = 800; J = 847.5; K = 2400; L = 2542.5; M = -2002; N = -3338.5 }
Set-Row "CUL" 58  @{ H = 3040; J = 3040; L = 9120; N = -9376 }
Set-Row "CUL" 93  @{ H = 4696.933; J = 4696.933; L = 14090.799; N = -17834.799 }
Set-Row "CUL" 98  @{ H = 492.4; I = 400; J = 605.3333; K = 1200; L = 1815.9999; M = 298; N = -4811.9999 }
Set-Row "CUL" 110 @{ H = 3293.7144; I = 0; J = 3293.7144; K = 0; L = 9881.143199999999; M = $null; N = -18061.1432 }

# ---- LTW ----
Set-Row "LTW" 68  @{ H = 0; I = 0; J = 0; K = 0; L = 0; M = $null; N = $null }
Set-Row "LTW" 71  @{ H = 0; I = 0; J = 0; K = 0; L = 0; M = $null; N = $null }

# ---- WVR ----
Set-Row "WVR" 132 @{ H = 1695.2941; I = 622.1; J = 3228.4285; K = 1866.3; L = 9685.2855; M = 663.6999999999998; N = -14745.2855 }
